$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.163.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.06%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.477.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.15%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.42%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  +1.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.475.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.11%  "

$ws.Range("E10").Value = "  +0.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.166"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.71%  "

$ws.Range("E12").Value = "  -2.09%  "

$ws.Range("E13").Value = "  +1.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "69.022.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.23%  "

$ws.Range("E15").Value = "  +0.38%  "

$ws.Range("E16").Value = "  -0.60%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "23.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.514.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.42%  "

$ws.Range("E19").Value = "  -1.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "339.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.76%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.28%  "

$ws.Range("E22").Value = "  +0.42%  "

$ws.Range("E23").Value = "  +3.11%  "

$ws.Range("E24").Value = "  +0.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.65%  "

$ws.Range("E26").Value = "  -0.97%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.604.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.54%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.46%  "

$ws.Range("E30").Value = "  -1.63%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "434.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.73%  "

$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("E34").Value = "  -1.60%  "

$ws.Range("E35").Value = "  -2.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.19%  "

$ws.Range("E37").Value = "  +0.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.111"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.52%  "

$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.23%  "

$ws.Range("E41").Value = "  -1.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.89%  "

$ws.Range("E43").Value = "  -3.05%  "

$ws.Range("E44").Value = "  +1.37%  "

$ws.Range("E45").Value = "  +0.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.41%  "

$ws.Range("E47").Value = "  +0.38%  "

$ws.Range("E48").Value = "  +0.37%  "

$ws.Range("E49").Value = "  -0.09%  "

$ws.Range("E50").Value = "  +1.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0916"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.25%  "
